$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date of the report
$ws.Range("B3").Value = [DateTime]"2021-10-07"

# Team name / numeric values entered by the team
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 5

# Team member names (replace generic placeholders with real names)
$ws.Range("A8").Value = "Lothaire Aubergeon"
$ws.Range("A9").Value = "Ian Hutter"
$ws.Range("A10").Value = "Sophie Leichtle"
$ws.Range("A11").Value = "Marco Schöb"
$ws.Range("A12").Value = "Markus Wagner"
$ws.Range("B12").Value = "4.5h"

# Tasks completed / next week
$ws.Range("A19").Value = """How might we"" statements formulated"
$ws.Range("B19").Value = "what comes next by the professor :)"
$ws.Range("A20").Value = "looked for sources"
$ws.Range("A21").Value = "summarized sources"
$ws.Range("A22").Value = "brainstorming"

# Active cell moved to B13 after editing
$ws.Range("B13").Select()

# Printer/page setup info saved with the workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
